$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: add the new TC id in column A and change the BUSQUEDA value (col D) ---
# (value order matters for the shared-strings table: CVBNM is written before DEC_0711)
$ws.Range("D3").Value = "CVBNM"
$ws.Range("A3").Value = "DEC_0711"

# --- Row 4: the row previously only carried empty, pre-formatted B4/C4 cells; ---
# --- now it becomes a full data row for the new TC script DEC_0715 ---
$ws.Range("A4").Value = "DEC_0715"
$ws.Range("B4").Value = "18092588-0"

# C4 already carried a "quote-prefix" text style (s=3); a plain value write resets
# that style, so restore it by pasting formats only from the sibling cell C3.
$ws.Range("C4").Value = "sebA`$1357"
$ws.Range("C3").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D4").Value = "SIN_DATO"
$ws.Range("E4").Value = "SIN_DATO"
$ws.Range("F4").Value = "SIN_DATO"
$ws.Range("G4").Value = "SIN_DATO"
$ws.Range("H4").Value = "SIN_DATO"
$ws.Range("I4").Value = "SIN_DATO"
$ws.Range("J4").Value = "SIN_DATO"

# --- Column D width change (stored width 10.44140625 -> 14) ---
$ws.Columns("D").ColumnWidth = 13.166666666666666

# --- Selection moves to A5 ---
$ws.Range("A5").Select()
